$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift the content of rows 11-21 down to rows 12-22 (working bottom-up so we don't
# clobber a value before it has been copied down). Using Range.Copy preserves both the
# original cell type (string vs boolean, e.g. the literal text "true") and the cell
# style, unlike a plain .Value assignment which can coerce "true"/"false" text into a
# real Boolean. The destination is cleared first since Copy() leaves a stale value in
# place when the source cell is blank.
for ($r = 21; $r -ge 11; $r--) {
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $dst.ClearContents()
    $ws.Range("A" + $r + ":B" + $r).Copy($dst)
}

# Insert the new "Jurisdiction" property row at row 11 (reuses row 11's existing style,
# which is now a copy of the original row 11 style since nothing but values changed)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "0.1.1"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"
